# Updated cryptos list (Price/Volume(1h) refresh, plus a few reordered rows)
# matching the GitHub Actions scraper commit.
#
# Note: Price values in column D are prefixed with a literal leading
# apostrophe ('). This forces Excel to store them as plain text, exactly
# as they appeared in the source data (e.g. "0.0000133" or "18.70"),
# instead of auto-converting them into numbers (which would turn
# "0.0000133" into 1.33E-05 or drop the trailing zero from "18.70").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''57.698.18'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '''2.565.81'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''515.42'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").Value = '''141.74'
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  -1.56%  '
$ws.Range("D9").Value = '''2.581.30'
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").Value = '''6.65'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '''0.101'
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D12").Value = '''0.322'
$ws.Range("E12").Value = '  -4.65%  '
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").Value = '''3.019.88'
$ws.Range("E14").Value = '  -2.93%  '
$ws.Range("D15").Value = '''57.720.51'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("E16").Value = '  -3.36%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000133'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '''2.575.41'
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").Value = '''334.73'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").Value = '''10.17'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").Value = '''6.29'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '''65.15'
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  -6.15%  '
$ws.Range("D28").Value = '''2.682.51'
$ws.Range("E28").Value = '  -3.08%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '''0.0₃0738'
$ws.Range("E31").Value = '  -7.26%  '
$ws.Range("D32").Value = '''6.19'
$ws.Range("E32").Value = '  -7.01%  '
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '''149.93'
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''18.64'
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("E36").Value = '  -3.96%  '
$ws.Range("D37").Value = '''1.13'
$ws.Range("E37").Value = '  -4.20%  '
$ws.Range("D38").Value = '''0.866'
$ws.Range("E38").Value = '  -4.23%  '
$ws.Range("D39").Value = '''35.88'
$ws.Range("E39").Value = '  -2.66%  '
$ws.Range("E40").Value = '  -3.74%  '
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '''268.63'
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '''10.65'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = '''0.0947'
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").Value = '''18.70'
$ws.Range("E48").Value = '  -3.60%  '
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("D50").Value = '''1.968.16'
$ws.Range("E50").Value = '  -3.91%  '
$ws.Range("D51").Value = '''4.59'
$ws.Range("E51").Value = '  -1.93%  '
